$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column G ("Notes") content updates ---
# Row 3: text unchanged (kept for completeness / no-op safe to set)
$ws.Range("G3").Value = "Reason: No matching Account Number on the lookup table. Account Number: TESTING"

# Rows 4-7: clear the stale "Reason: ..." notes
$ws.Range("G4").ClearContents()
$ws.Range("G5").ClearContents()
$ws.Range("G6").ClearContents()
$ws.Range("G7").ClearContents()

# Rows 8-14: replace old notes with the new "AppliedPatch_*" labels
$ws.Range("G8").Value = "AppliedPatch_OnlyQueue"
$ws.Range("G9").Value = "AppliedPatch_CCR"
$ws.Range("G10").Value = "AppliedPatch_AssignTo"
$ws.Range("G11").Value = "AppliedPatch_CCRAssignTo"
$ws.Range("G12").Value = "AppliedPatch_CCRAssignTo"
$ws.Range("G13").Value = "AppliedPatch_CCR"
$ws.Range("G14").Value = "AppliedPatch_AssignTo"

# --- Row heights (wrap-text rows re-measured at new font metrics) ---
$ws.Rows.Item(2).RowHeight = 45
$ws.Rows.Item(3).RowHeight = 45
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(5).RowHeight = 30
$ws.Rows.Item(6).RowHeight = 45
$ws.Rows.Item(7).RowHeight = 45
$ws.Rows.Item(8).RowHeight = 30
$ws.Rows.Item(9).RowHeight = 30
$ws.Rows.Item(10).RowHeight = 30
$ws.Rows.Item(11).RowHeight = 45

# --- Column A width widened slightly, bestFit turned off ---
# (19.2 chars -> engine stores padded width 20, matching the target XML)
$ws.Columns.Item(1).ColumnWidth = 19.2

# --- Selection moved from E17 to F11, no more frozen/scrolled topLeftCell ---
$ws.Range("F11").Select()
